$d = $word.ActiveDocument

# The document currently ends with (1-based paragraph indices):
#   18: I've implemented this intermediate step as agent2.py.
#   19: <empty paragraph>
#   20: <empty paragraph that carries the _GoBack bookmark>
#
# We need to insert a "Part 3" section before the bookmarked paragraph and
# then add the final body text to that paragraph itself, just ahead of the
# bookmark it carries. Always re-fetch paragraphs by fresh index - cached
# paragraph object references do not reliably track subsequent edits in
# this host.

$bookmarkIndex = 20

# Insert 7 new (empty) paragraphs immediately before the bookmarked
# paragraph. Each InsertParagraphBefore() call pushes the bookmarked
# paragraph one slot further down, so we re-target the same (shifting)
# index each time.
for ($i = 0; $i -lt 7; $i++) {
    $d.Paragraphs.Item($bookmarkIndex).Range.InsertParagraphBefore()
    $bookmarkIndex = $bookmarkIndex + 1
}

# After the loop the paragraphs are laid out like this:
#   20: blank
#   21: "Part 3: Implement Q-Learning"  (bold heading)
#   22: blank
#   23: Q-learning implementation paragraph
#   24: blank
#   25: deadlines / trials paragraph
#   26: blank
#   27: (originally 20) bookmarked paragraph -> gets the closing body text

$d.Paragraphs.Item(21).Range.Text = "Part 3: Implement Q-Learning"
$d.Paragraphs.Item(21).Range.Font.Bold = 1

$d.Paragraphs.Item(23).Range.Text = "I implemented Q-Learning with the learning _rate = 0.5 and discount_factor = 0.5.  The Q dictionary is initialized with all values set to 0 by default.  At every step, the agent checks the Q-values for all possible actions (None, left, right, and forward) and picks the one that yields the largest Q-value.  If all possible actions lead to a Q-value of 0 (which means those Q-values have never been calculated yet), then the agent will just pick one of the 4 actions randomly."

$d.Paragraphs.Item(25).Range.Text = "I also started to enforce deadlines and changed the # of trials to 100."

$apostrophe = [string][char]8217
$closingText = "With the above implementation, I found that my agent would quickly learn not to disobey the traffic signals after several trials.  However, it would get stuck in a local optima where it ends up staying in the same place (action = None) and collecting a reward of 1 every time.  Since the traffic lights seem to cycle every 5 rounds, I" + $apostrophe + "ve modified my agent to keep track of how many rounds in a row did it choose action = None and if it" + $apostrophe + "s greater than 6, then it chooses an action at random regardless of Q values.  With this change, the agent doesn" + $apostrophe + "t get stuck in one place and can continue learning."

$d.Paragraphs.Item($bookmarkIndex).Range.InsertBefore($closingText)
